# Slide 7 ("Conclusões"), body placeholder (Shape 2 / id 171):
#  - shrink the bullet text (sz 31.23 -> 26pt)
#  - re-wrap / extend bullet 2's wording (split across extra runs)
#  - replace bullet 3's wording and fold the old blank line + bullet 4
#    ("Já começámos...") into it, then delete the now-unused paragraphs

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Bullet 1: "Podemos afirmar..." -> just a font-size change ---
$para1 = $tr.Paragraphs(1)
$para1.Font.Size = 26

# --- Bullet 2: "Foram planeados..." -> re-split into 4 runs, smaller font ---
$para2 = $tr.Paragraphs(2)
$para2.Text = "Foram planeados e analisados os vários elementos do sistema, a forma como eles devem ser "
$run2b = $para2.InsertAfter("construídos, quais ")
$run2c = $para2.InsertAfter("as suas ")
$run2d = $para2.InsertAfter("funções e as tecnologias e conceitos utilizados na sua implementação. ")
$para2.Font.Size = 26

# --- Bullet 3: "Debatemos..." -> replaced with the new progress text, ---
# --- absorbing the old blank paragraph + "Já começámos..." paragraph  ---
$para3 = $tr.Paragraphs(3)
$para3.Text = "O desenvolvimento dos elementos do sistema já se encontra perto da sua conclusão, restando apenas finalizar os testes, finalizar alguns ecrãs da aplicação móvel e melhorar o OCR"
$run3b = $para3.InsertAfter(".")
$para3.Font.Size = 26
$run3b.Font.Size = 31.23

# The paragraph that used to hold the blank line (now still index 4) and
# the one that used to hold "Já começámos..." (now shifted to index 4
# after the first delete) are no longer needed.
$tr.Paragraphs(4).Delete()
$tr.Paragraphs(4).Delete()
